# UP-Sell and Customer Churn test cases added
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # up_sell_Filter
$ws2 = $wb.Worksheets.Item(2)   # up_Sell_Report_EtoE

# ---------------------------------------------------------------------------
# Sheet 1: up_sell_Filter
# ---------------------------------------------------------------------------

# Row 2: fill in the Segmentation value ("Overall"), styled like the existing
# hyperlink-ish cells (B3 already carries this style).
$ws1.Range("B2").Value = "Overall"
$ws1.Range("B2").Style = "Hyperlink"

# New row 3: a second filter row for "Deposit Account" / "Within 90 days",
# matching the layout/format of row 2's Propensity & Date Range cells.
$ws1.Range("A3").Value = 2
$ws1.Range("C3").Value = "Deposit Account"
$ws1.Range("D3").Value = "Within 90 days"
$ws1.Range("C2:D2").Copy()
$ws1.Range("C3:D3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet 2: up_Sell_Report_EtoE
# ---------------------------------------------------------------------------

# Insert a new "Selected Drivers" column before the existing "Drivers Title"
# column (old D shifts to E, old E shifts to F).
$ws2.Columns("D").Insert()
$ws2.Range("D1").Value = "Selected Drivers"

$ws2.Range("D2").Value = "Driver 1,Driver 2,Driver 3,Driver 4,Driver 5"
$ws2.Range("D3").Value = "NA"
$ws2.Range("D4").Value = "Driver 1,Driver 2,Driver 3,Driver 4,Driver 5"
$ws2.Range("D5").Value = "NA"
$ws2.Range("D6").Value = "NA"
$ws2.Range("D7").Value = "Driver 1,Driver 2,Driver 3,Driver 4,Driver 5"

# Rows 4 and 7 swap which Customer-Probability/Propensity combination they
# report (row 4 becomes the "100%-90% / Extreme Likely" case, row 7 becomes
# the "100%-50% / Extreme Likely,High Likely,Likely" case).
$ws2.Range("B4").Value = "100%-90%"
$ws2.Range("C4").Value = "Extreme Likely"
$ws2.Range("E4").Value = "Product Up-Sell Drivers for Top 90% to 100%"
$ws2.Range("F4").Value = "Profiles for Top 90% to 100%"

$ws2.Range("B7").Value = "100%-50%"
$ws2.Range("C7").Value = "Extreme Likely,High Likely,Likely"
$ws2.Range("E7").Value = "Product Up-Sell Drivers for Top 50% to 100%"
$ws2.Range("F7").Value = "Profiles for Top 50% to 100%"

# ---------------------------------------------------------------------------
# Selections (sheet1 first, sheet2 last so sheet2 stays the active tab)
# ---------------------------------------------------------------------------
$ws1.Range("C11").Select()
$ws2.Range("E12").Select()
